$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to text format so numeric-looking price strings are not
# auto-converted to floating point numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '64.922.20'
$ws.Range("E2").Value = '  +0.24%  '

$ws.Range("D3").Value = '3.157.49'
$ws.Range("E3").Value = '  -0.03%  '

$ws.Range("E4").Value = '  +0.03%  '

$ws.Range("D5").Value = '579.12'
$ws.Range("E5").Value = '  +1.12%  '

$ws.Range("D6").Value = '149.92'
$ws.Range("E6").Value = '  -1.00%  '

$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  +0.03%  '

$ws.Range("D8").Value = '3.153.27'
$ws.Range("E8").Value = '  -0.07%  '

$ws.Range("D9").Value = '0.527'
$ws.Range("E9").Value = '  -0.39%  '

$ws.Range("D10").Value = '0.159'
$ws.Range("E10").Value = '  -2.33%  '

$ws.Range("D11").Value = '6.11'
$ws.Range("E11").Value = '  -0.41%  '

$ws.Range("D12").Value = '0.502'
$ws.Range("E12").Value = '  -0.79%  '

$ws.Range("D13").Value = '0.0000266'
$ws.Range("E13").Value = '  +3.29%  '

$ws.Range("D14").Value = '37.33'
$ws.Range("E14").Value = '  -2.32%  '

$ws.Range("D15").Value = '3.678.38'
$ws.Range("E15").Value = '  +0.04%  '

$ws.Range("D16").Value = '64.910.37'
$ws.Range("E16").Value = '  +0.06%  '

$ws.Range("B17").Value = 'WrappedEther'
$ws.Range("C17").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D17").Value = '3.161.40'
$ws.Range("E17").Value = '  +0.14%  '

$ws.Range("B18").Value = 'Polkadot'
$ws.Range("C18").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D18").Value = '7.16'
$ws.Range("E18").Value = '  -1.22%  '

$ws.Range("D19").Value = '0.111'
$ws.Range("E19").Value = '  +0.44%  '

$ws.Range("D20").Value = '505.97'
$ws.Range("E20").Value = '  -2.60%  '

$ws.Range("D21").Value = '14.96'
$ws.Range("E21").Value = '  -0.35%  '

$ws.Range("D22").Value = '0.718'
$ws.Range("E22").Value = '  -2.81%  '

$ws.Range("D23").Value = '15.28'
$ws.Range("E23").Value = '  +0.27%  '

$ws.Range("D24").Value = '7.74'
$ws.Range("E24").Value = '  -1.56%  '

$ws.Range("D25").Value = '84.61'
$ws.Range("E25").Value = '  -0.91%  '

$ws.Range("D26").Value = '0.999'
$ws.Range("E26").Value = '  -0.01%  '

$ws.Range("D27").Value = '9.00'
$ws.Range("E27").Value = '  +2.30%  '

$ws.Range("E28").Value = '  -0.31%  '

$ws.Range("D29").Value = '2.19'
$ws.Range("E29").Value = '  -0.20%  '

$ws.Range("D30").Value = '2.80'
$ws.Range("E30").Value = '  +4.47%  '

$ws.Range("D31").Value = '27.67'
$ws.Range("E31").Value = '  -1.30%  '

$ws.Range("E32").Value = '  +0.01%  '

$ws.Range("E33").Value = '  +1.09%  '

$ws.Range("D34").Value = '6.24'
$ws.Range("E34").Value = '  +1.82%  '

$ws.Range("D35").Value = '6.50'
$ws.Range("E35").Value = '  -1.46%  '

$ws.Range("D36").Value = '54.88'
$ws.Range("E36").Value = '  -1.70%  '

$ws.Range("B37").Value = 'Bittensor'
$ws.Range("C37").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D37").Value = '486.03'
$ws.Range("E37").Value = '  -0.96%  '

$ws.Range("B38").Value = 'Hedera'
$ws.Range("C38").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D38").Value = '0.0892'
$ws.Range("E38").Value = '  +2.73%  '

$ws.Range("D39").Value = '0.0418'
$ws.Range("E39").Value = '  -1.53%  '

$ws.Range("D40").Value = '2.95'
$ws.Range("E40").Value = '  -1.39%  '

$ws.Range("D41").Value = '8.79'
$ws.Range("E41").Value = '  +1.18%  '

$ws.Range("D42").Value = '3.003.98'
$ws.Range("E42").Value = '  -3.66%  '

$ws.Range("D43").Value = '0.115'
$ws.Range("E43").Value = '  -4.15%  '

$ws.Range("D44").Value = '2.44'
$ws.Range("E44").Value = '  -1.04%  '

$ws.Range("D45").Value = '0.283'
$ws.Range("E45").Value = '  -4.87%  '

$ws.Range("D46").Value = '28.35'
$ws.Range("E46").Value = '  -3.32%  '

$ws.Range("D47").Value = '0.0₃0592'
$ws.Range("E47").Value = '  +2.04%  '

$ws.Range("D48").Value = '1.00'

$ws.Range("E49").Value = '  -1.74%  '

$ws.Range("D50").Value = '2.25'
$ws.Range("E50").Value = '  -2.73%  '

$ws.Range("D51").Value = '2.50'
$ws.Range("E51").Value = '  +15.67%  '

# Restore the original (default/general) style for column D now that the
# text values are safely stored, so the cell style stays unchanged.
$ws.Range("D2:D51").Style = "Normal"
